$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new employee row above row 16 (CARLOS ANDRES ACOSTA TERAN), shifting
# the rest of the table (and the footer signature rows) down by one.
$ws.Rows.Item(16).Insert()

# The freshly inserted row has no formatting - clone it from the row right
# below (which still carries the original table-row formatting) so borders /
# fonts / fills / number formats match the rest of the table.
$ws.Range("B17:J17").Copy()
$ws.Range("B16:J16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New employee: FERNANDO JOSE GUEVARA CANCHILA
$ws.Cells.Item(16, 2).Value = "CC"
$ws.Cells.Item(16, 3).Value = "92523071"
$ws.Cells.Item(16, 4).Value = "FERNANDO JOSE GUEVARA CANCHILA"
$ws.Cells.Item(16, 5).Value = "2509"
$ws.Cells.Item(16, 6).Value = 56940
$ws.Cells.Item(16, 7).Value = 1423500

# Update "Periodo Mora" for every other worker row (2508 -> 2509)
$ws.Cells.Item(17, 5).Value = "2509"
$ws.Cells.Item(18, 5).Value = "2509"
$ws.Cells.Item(19, 5).Value = "2509"
$ws.Cells.Item(20, 5).Value = "2509"
$ws.Cells.Item(21, 5).Value = "2509"

# OSCAR LUIS TORRES RODRIGUEZ (row 20) now has a partial-period "Valor Mora"
$ws.Cells.Item(20, 6).Value = 30368

# Header summary figures
$ws.Range("E11").Value = 321288
$ws.Range("C13").Value = 6

$excel.ActiveWorkbook.Save()
